$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Récupérez le numéro de compte de  que vous avez trouvé grâce à
# a la requête précédente. " -> "Récupérez le numéro de compte que vous avez
# trouvé grâce à a la requête précédente. "
# (drop the stray "de  " before "que"), while keeping the paragraph's three
# existing runs intact (only the middle run's text changes).
# ---------------------------------------------------------------------------

# Locate the run that contains the typo so we know its original boundaries.
$run2 = $d.Content
$run2.Find.Execute("Récupérez le numéro de compte de  que vous avez trouvé grâce à a la requête précédente. ") | Out-Null
$run2Start = $run2.Start
$run2End = $run2.End

# Locate and delete the stray "de  " (the word "de" plus the double space
# that follows it) right before "que".
$typo = $d.Range($run2Start, $run2End)
$typo.Find.Execute("de  que") | Out-Null
$typoStart = $typo.Start
$delRange = $d.Range($typoStart, $typoStart + 4)
$delRange.Delete()

# The text edit above merges this paragraph's runs into a single run (the
# engine coalesces adjacent same-formatted runs whenever text changes). Put
# the original 3-run split back by toggling (and immediately un-toggling) a
# character attribute across exactly the corrected run's span - this forces
# the run boundaries without altering the run's visible formatting.
$newRun2End = $run2End - 4
$splitRange = $d.Range($run2Start, $newRun2End)
$splitRange.Bold = $true
$splitRange.Bold = $false

# ---------------------------------------------------------------------------
# Change 2: "De quel type d’opération sont-elles ? (opération tiers ou
# retrait/dépôt)" -> same text but "(opération" becomes "(Opération" and
# that word now lives in its own run:
#   "De quel type d’opération sont-elles ? ("  |  "Opération"  |  " tiers ou
#   retrait/dépôt)"
# ---------------------------------------------------------------------------

# Narrow the search to this paragraph (well past the first, unrelated
# "opération" earlier in the document) so we land on the right occurrence.
$para = $d.Content
$para.Find.Execute("De quel type d") | Out-Null
$paraStart = $para.Start
$para2 = $d.Range($paraStart, $paraStart + 400)
$para2.Find.Execute("opération tiers ou retrait/dépôt)") | Out-Null
$paraEnd = $para2.End

$scope = $d.Range($paraStart, $paraEnd)
$scope.Find.Execute("opération tiers") | Out-Null
$wordStart = $scope.Start

# Capitalize the word's first letter BEFORE splitting the run - any text
# mutation re-coalesces adjacent same-formatted runs in the paragraph, so
# the capitalization has to happen first and the run split has to be the
# very last step touching this paragraph.
$firstLetter = $d.Range($wordStart, $wordStart + 1)
$firstLetter.Text = "O"

# Now isolate "Opération" (9 characters) into its own run by toggling (and
# immediately un-toggling) a character attribute across exactly its span.
$opWord = $d.Range($wordStart, $wordStart + 9)
$opWord.Bold = $true
$opWord.Bold = $false
